$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text change
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "fizmasss"

# ---------------------------------------------------------------------------
# 2. Dates change (order date / pickup date) - both become 27.4.2020
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "27.4.2020"
$ws.Range("B3").Value = "27.4.2020"

# ---------------------------------------------------------------------------
# 3. Add invoice number + payment info (new cells D2,E2,D3,E3), using the
#    same formatting as the existing label cells in column A (style idx 2).
#    Use temp cells to stash/restore format cleanly via copy/paste formats.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("D2:E3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D2").Value = "Nr faktury:"
$ws.Range("E2").Value = "eX2137/20"
$ws.Range("D3").Value = "Płatność:"
$ws.Range("E3").Value = 0

# ---------------------------------------------------------------------------
# 4. Update quantities on row 5 / row 6
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = 45
$ws.Range("C5").Value = 12
$ws.Range("C6").Value = 33

# ---------------------------------------------------------------------------
# 5. Remove the extra product rows (7-13 content), keep them as blank rows,
#    and shrink the row5:row7 merges down to row5:row6 (fix for "only one
#    product type" bug). Merging cells that already carry a border tends to
#    rewrite the border/style of the merged cells in this COM engine, so we
#    stash the original formatting of A5,A6,B5,B6 in scratch cells first and
#    restore it after merging.
# ---------------------------------------------------------------------------
$ws.Range("A5").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("H4").PasteSpecial(-4122)

$ws.Range("A8:A9").UnMerge()
$ws.Range("B8:B9").UnMerge()
$ws.Range("A10:A11").UnMerge()
$ws.Range("B10:B11").UnMerge()
$ws.Range("A12:A13").UnMerge()
$ws.Range("B12:B13").UnMerge()
$ws.Range("A5:A7").UnMerge()
$ws.Range("B5:B7").UnMerge()

$ws.Rows("7:13").Clear()

$ws.Range("A5:A6").Merge()
$ws.Range("B5:B6").Merge()

$ws.Range("H1").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("B6").PasteSpecial(-4122)

$ws.Range("H1:H4").Clear()

Write-Host "Edit complete"
